$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed cells keep their existing text representation
# (values like "27.699.07" or "  -0.55%  " must not be reinterpreted as numbers/dates).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.699.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4883"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3792"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07327"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9132"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.55"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07661"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.895.17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.483"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.611"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008770"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.742.39"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.121"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.129.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.156"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.868"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08912"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.202"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.222"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7664"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02038"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.532"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.096"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5467"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.976"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.885"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.516"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "112.07"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.64"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4777"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.639"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.19%  "
